# Apply the changes described in the commit "fixing example in supplementary materials"
# to the genomics schema_and_data workbook.

$wb = $excel.ActiveWorkbook

# Sheets (by index, 1-based, in workbook order):
#   1 = "!!_Table of contents"
#   2 = "!!_Schema"
#   3 = "!!Genes"
#   4 = "!!Transcripts"
$wsSchema = $wb.Worksheets.Item(2)

# --- Fix the Schema sheet row ordering ---------------------------------
# Swap the "symbol" and "location" attribute rows of the Gene class
# (rows 5 and 6) so that "symbol" is listed before "location".
$a5 = $wsSchema.Range("A5").Value2
$b5 = $wsSchema.Range("B5").Value2
$c5 = $wsSchema.Range("C5").Value2
$d5 = $wsSchema.Range("D5").Value2
$e5 = $wsSchema.Range("E5").Value2

$a6 = $wsSchema.Range("A6").Value2
$b6 = $wsSchema.Range("B6").Value2
$c6 = $wsSchema.Range("C6").Value2
$d6 = $wsSchema.Range("D6").Value2
$e6 = $wsSchema.Range("E6").Value2

$wsSchema.Range("A5").Value = $a6
$wsSchema.Range("B5").Value = $b6
$wsSchema.Range("C5").Value = $c6
$wsSchema.Range("D5").Value = $d6
$wsSchema.Range("E5").Value = $e6

$wsSchema.Range("A6").Value = $a5
$wsSchema.Range("B6").Value = $b5
$wsSchema.Range("C6").Value = $c5
$wsSchema.Range("D6").Value = $d5
$wsSchema.Range("E6").Value = $e5

# Swap the "gene" and "id" attribute rows of the Transcript class
# (rows 8 and 9) so that "id" is listed before "gene".
$a8 = $wsSchema.Range("A8").Value2
$b8 = $wsSchema.Range("B8").Value2
$c8 = $wsSchema.Range("C8").Value2
$d8 = $wsSchema.Range("D8").Value2
$e8 = $wsSchema.Range("E8").Value2

$a9 = $wsSchema.Range("A9").Value2
$b9 = $wsSchema.Range("B9").Value2
$c9 = $wsSchema.Range("C9").Value2
$d9 = $wsSchema.Range("D9").Value2
$e9 = $wsSchema.Range("E9").Value2

$wsSchema.Range("A8").Value = $a9
$wsSchema.Range("B8").Value = $b9
$wsSchema.Range("C8").Value = $c9
$wsSchema.Range("D8").Value = $d9
$wsSchema.Range("E8").Value = $e9

$wsSchema.Range("A9").Value = $a8
$wsSchema.Range("B9").Value = $b8
$wsSchema.Range("C9").Value = $c8
$wsSchema.Range("D9").Value = $d8
$wsSchema.Range("E9").Value = $e8

# Fix the "five_prime" attribute's Format to mark it as the primary,
# unique key of the Location class.
$wsSchema.Range("D13").Value = "PositiveInteger(primary=True, unique=True)"

# --- Make the Schema sheet the active tab -------------------------------
$wsSchema.Activate()

$wb.Save()
